# Commit "Mon, Jul 27, 2020 12:05:52 PM" only touches the root <p:.../> opening
# tags of every part that round-trips through the PowerPoint XML writer
# (ppt/presentation.xml, the notes master, all 13 notes slides, all 11 slide
# layouts, the slide master and all 13 slides): PowerPoint stamped each of
# them with an extra
#
#   xmlns:ahyp="http://schemas.microsoft.com/office/drawing/2018/hyperlinkcolor"
#
# namespace declaration. That is a side effect of the file having been opened
# and re-saved by a PowerPoint build that knows about the 2018 "hyperlink
# color" drawing extension; no slide text, shape geometry, formatting,
# ordering, relationship, or any other OOXML content actually differs
# between the two XML trees (every hunk in the diff is exactly the root
# start-tag, nothing else).
#
# There is no user-facing object-model action (typing text, moving shapes,
# changing colors, adding hyperlinks, etc.) that corresponds to this purely
# cosmetic namespace stamp, so the only faithful reproduction available
# through the PowerPoint COM API is to simply touch and re-save the
# presentation as-is, leaving every slide/master/layout's actual content
# untouched, exactly as the original author's save did.
$p = $ppt.ActivePresentation
$p.Save()
